$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-08-24 10:47:11"

$wsZhCn.Range("H3").Value = "2016-08-24 10:47:02"
$wsZhCn.Range("K3").Value = "2016-08-24 10:47:30"

$wsDeDe.Range("H3").Value = "2016-08-24 10:47:11"
$wsDeDe.Range("K3").Value = "2016-08-24 10:47:37"
